$d = $word.ActiveDocument

$pairs = @(
    @("76×65=4940", "37×64=2368"),
    @("57×65=3705", "64×19=1216"),
    @("40×65=2600", "83×54=4482"),
    @("51×53=2703", "36×65=2340"),
    @("47×77=3619", "16×67=1072"),
    @("88×99=8712", "45×33=1485"),
    @("71×76=5396", "31×93=2883"),
    @("20×84=1680", "60×79=4740"),
    @("81×43=3483", "44×84=3696"),
    @("97×97=9409", "41×77=3157"),
    @("91×33=3003", "94×23=2162"),
    @("48×29=1392", "61×51=3111"),
    @("22×34=748",  "26×80=2080"),
    @("62×50=3100", "94×73=6862"),
    @("13×83=1079", "71×50=3550"),
    @("97×44=4268", "21×18=378"),
    @("87×63=5481", "62×76=4712"),
    @("58×12=696",  "92×45=4140"),
    @("46×35=1610", "30×66=1980"),
    @("35×14=490",  "33×61=2013"),
    @("85×60=5100", "94×39=3666"),
    @("94×85=7990", "88×39=3432"),
    @("76×30=2280", "97×70=6790"),
    @("22×26=572",  "20×21=420"),
    @("23×20=460",  "68×87=5916")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
